$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data rows (Cluster name, Active cases), sorted alphabetically by cluster name.
$data = @(
    @('3153 Sacred Heart Community St Kilda Tier 1A', 21),
    @('3528 Ottoman Village Aged Care Broadmeadows', 25),
    @('3600 Belvedere Aged Care Noble Park', 24),
    @('3612 BlueCross Glengowrie', 34),
    @('3652 Regis Aged Care Dandenong North', 24),
    @('3684 Homestyle Aged Care Langford Grange Cranbourne East', 15),
    @('3824 Estia Health South Morang', 48),
    @('3980 Arcare Keysborough Aged Care Keysborough', 23),
    @('4518 Regis Aged Care Fawkner', 13),
    @('ACFS Port Logistics Altona', 15),
    @('Alfred Health The Alfred Hospital Melbourne', 13),
    @('Armstrong Creek School Armstrong Creek', 17),
    @('Aspect Autism Spectrum Australia Disability Service Heatherton', 12),
    @('Berwick Fields Primary School Berwick', 13),
    @('Berwick Lodge Primary School Berwick', 17),
    @('Bridgewood Primary School Officer', 10),
    @('Bubup Womindjeka Family and Children''s Centre Port Melbourne', 10),
    @('Clifton Hill Primary School Clifton Hill', 15),
    @('Community Kids Pascoe Vale Early Education Centre Pascoe Vale South', 11),
    @('Dandenong North Primary School', 10),
    @('Elements Childcare Warralily Armstrong Creek', 27),
    @('Highett St Apartment Complex Richmond', 14),
    @('Honeyeater Hairdressers Bendigo', 13),
    @('Inghams Chicken Warehouse Lyndhurst', 10),
    @('Inghams Enterprises Somerville', 18),
    @('JBS Australia Brooklyn', 14),
    @('Kmart Distribution Centre Truganina', 12),
    @('Lighthouse Christian College', 13),
    @('Lilydale Motor Inn Lilydale', 10),
    @('Lowanna College Newborough', 15),
    @('McQuinns Gym Bendigo', 24),
    @('Melbourne Alternative Place of Detention Carlton', 33),
    @('Metcash Limited Distribution Centre Laverton North', 19),
    @('Monash Health Casey Hospital Emergency Department Tier 1B', 15),
    @('Nido Early School Wyndham Vale', 13),
    @('Pelican Childcare Craigieburn', 13),
    @('Shawlands Caravan Park Dandenong South', 14),
    @('St Mary''s Primary School Swan Hill', 14),
    @('St Vincents Hospital Emergency Department Melbourne', 26),
    @('TUROSI PTY LTD Thomastown', 13),
    @('The Robin Hood Inn Drouin West', 11),
    @('The Royal Children''s Hospital Parkville', 11),
    @('Vizzarri Farms Koo Wee Rup', 26),
    @('Werribee Mercy Hospital Emergency Department', 43),
    @('Western Health Sunshine Hospital Emergency Department', 22)
)

# Clear any existing data rows below the header first (old sheet had data through row 44).
$oldUsedRange = $ws.UsedRange
$lastRow = $oldUsedRange.Rows.Count
if ($lastRow -gt 1) {
    $ws.Range("A2:B$lastRow").Clear()
}

$rowCount = $data.Count
$startRow = 2
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
